$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures per the Jan 26 2023 GitHub Actions data refresh.
# Values are plain text (matching the sheet's existing inline-string cells), so they are
# entered with a leading apostrophe to stop Excel from auto-converting them to numbers/percents.
$ws.Range("D2").Value = "'307.64"
$ws.Range("E2").Value = "'2.29%"
$ws.Range("D3").Value = "'35.92"
$ws.Range("E3").Value = "'1.38%"
$ws.Range("D4").Value = "'5.044"
$ws.Range("E4").Value = "'-0.07%"
$ws.Range("D5").Value = "'0.08131"
$ws.Range("E5").Value = "'1.88%"
$ws.Range("D6").Value = "'1.954"
$ws.Range("E6").Value = "'2.28%"
$ws.Range("D7").Value = "'4.132"
$ws.Range("E7").Value = "'1.59%"
$ws.Range("D8").Value = "'7.793"
$ws.Range("E9").Value = "'0.90%"
$ws.Range("D10").Value = "'0.1332"
$ws.Range("E10").Value = "'-6.02%"
$ws.Range("D11").Value = "'0.1918"
$ws.Range("E11").Value = "'0.68%"
$ws.Range("D12").Value = "'0.09239"
$ws.Range("E12").Value = "'0.29%"
$ws.Range("D13").Value = "'0.03506"
$ws.Range("E13").Value = "'2.88%"
$ws.Range("D14").Value = "'0.09868"
$ws.Range("E14").Value = "'-0.09%"
$ws.Range("D15").Value = "'0.001414"
$ws.Range("E15").Value = "'1.96%"
$ws.Range("D16").Value = "'0.005794"
$ws.Range("E16").Value = "'-0.06%"
$ws.Range("E17").Value = "'2.39%"
$ws.Range("E18").Value = "'-1.38%"
$ws.Range("D19").Value = "'0.3430"
$ws.Range("E19").Value = "'0.82%"
$ws.Range("E20").Value = "'3.64%"
$ws.Range("D21").Value = "'5.189"
$ws.Range("E21").Value = "'2.64%"
$ws.Range("E22").Value = "'7.86%"
$ws.Range("D23").Value = "'0.04379"
$ws.Range("E23").Value = "'-2.76%"
$ws.Range("E24").Value = "'0.35%"
$ws.Range("D25").Value = "'0.004770"
$ws.Range("E25").Value = "'-0.31%"
$ws.Range("E26").Value = "'33.88%"
$ws.Range("E27").Value = "'3.86%"
$ws.Range("D39").Value = "'0.01998"
$ws.Range("E39").Value = "'4.65%"
$ws.Range("D40").Value = "'0.05063"
$ws.Range("E40").Value = "'6.92%"
$ws.Range("D41").Value = "'0.01119"
$ws.Range("E41").Value = "'15.79%"
$ws.Range("D42").Value = "'0.007608"
$ws.Range("E42").Value = "'3.63%"
$ws.Range("E43").Value = "'3.95%"
$ws.Range("E44").Value = "'-0.91%"
$ws.Range("D45").Value = "'0.01128"
$ws.Range("E45").Value = "'8.74%"
$ws.Range("D46").Value = "'0.00006381"
$ws.Range("E46").Value = "'1.87%"
$ws.Range("E47").Value = "'-0.40%"
$ws.Range("D48").Value = "'63.57"
$ws.Range("E48").Value = "'-1.41%"
$ws.Range("E49").Value = "'-28.53%"
$ws.Range("E50").Value = "'-0.40%"
$ws.Range("E51").Value = "'-0.40%"
